# ===========================================================================
# Apply "Holden" simulation scheme update:
#  - insert 4 new data rows (Holden2.5 / Holden5 / Holden10 / Holden15)
#  - reorder + refresh the [h,k,l] columns (C:J) and refresh all other
#    computed statistic columns to the values produced by the updated run
# ===========================================================================
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Step 1: insert 4 blank rows before the HexGrid block (originally row 16) ----
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()
$ws.Rows.Item(16).Insert()

# ---- Step 2: match the index-column style (bold/centered/bordered) used elsewhere ----
$ws.Range("A15").Copy()
$ws.Range("A16:A19").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---- Step 3: update the [h,k,l] header labels (row 2, columns C:J) ----
$ws.Range("C2").Value = "[2, 1, 1]"
$ws.Range("D2").Value = "[4, 0, 0]"
$ws.Range("E2").Value = "[2, 0, 0]"
$ws.Range("F2").Value = "[2, 2, 0]"
$ws.Range("G2").Value = "[1, 1, 0]"
$ws.Range("H2").Value = "[3, 1, 0]"
$ws.Range("I2").Value = "[2, 2, 2]"
$ws.Range("J2").Value = "[3, 2, 1]"

# ---- Step 4: refresh data rows 3-15 (BT8Hex / Spiral / Offset schemes) ----
# row 3: BT8Hex_2.5
$ws.Range("C3").Value = 0.989639002269308
$ws.Range("D3").Value = 0.9824975127660879
$ws.Range("E3").Value = 0.9824975127660879
$ws.Range("F3").Value = 1.050541219657053
$ws.Range("G3").Value = 1.050541219657053
$ws.Range("H3").Value = 0.9555350363648139
$ws.Range("I3").Value = 0.9816293827479792
$ws.Range("J3").Value = 1.002537094226045
$ws.Range("K3").Value = 1.050541219657053
$ws.Range("L3").Value = 0.989639002269308
$ws.Range("M3").Value = 0.986068257517698
$ws.Range("N3").Value = 0.986068257517698
$ws.Range("O3").Value = 0.9758905171334034
$ws.Range("P3").Value = 1.007559244897483
$ws.Range("Q3").Value = 1.007559244897483
$ws.Range("R3").Value = 1.018304738587376
$ws.Range("S3").Value = 1.018304738587376
$ws.Range("T3").Value = 0.9937298746718811

# row 4: BT8Hex_5
$ws.Range("C4").Value = 0.9540881241053801
$ws.Range("D4").Value = 0.8833470588205872
$ws.Range("E4").Value = 0.8833470588205872
$ws.Range("F4").Value = 1.097522422524745
$ws.Range("G4").Value = 1.097522422524745
$ws.Range("H4").Value = 0.9090543681114364
$ws.Range("I4").Value = 0.9505537298916452
$ws.Range("J4").Value = 0.9847854886103488
$ws.Range("K4").Value = 1.097522422524745
$ws.Range("L4").Value = 0.9540881241053801
$ws.Range("M4").Value = 0.9187175914629837
$ws.Range("N4").Value = 0.9187175914629837
$ws.Range("O4").Value = 0.9154965170124679
$ws.Range("P4").Value = 0.978319201816904
$ws.Range("Q4").Value = 0.978319201816904
$ws.Range("R4").Value = 1.008120006993864
$ws.Range("S4").Value = 1.008120006993864
$ws.Range("T4").Value = 0.9632251986773571

# row 5: BT8Hex_10
$ws.Range("C5").Value = 0.9991217493748569
$ws.Range("D5").Value = 0.6192801061974191
$ws.Range("E5").Value = 0.6192801061974191
$ws.Range("F5").Value = 0.9944787308125864
$ws.Range("G5").Value = 0.9944787308125864
$ws.Range("H5").Value = 0.9251059867396035
$ws.Range("I5").Value = 0.6951985739839682
$ws.Range("J5").Value = 0.8731307857573206
$ws.Range("K5").Value = 0.9944787308125864
$ws.Range("L5").Value = 0.9991217493748569
$ws.Range("M5").Value = 0.8092009277861381
$ws.Range("N5").Value = 0.8092009277861381
$ws.Range("O5").Value = 0.8478359474372933
$ws.Range("P5").Value = 0.8709601954616208
$ws.Range("Q5").Value = 0.8709601954616208
$ws.Range("R5").Value = 0.9018398292993622
$ws.Range("S5").Value = 0.9018398292993622
$ws.Range("T5").Value = 0.8510526554776258

# row 6: BT8Hex_15
$ws.Range("C6").Value = 0.8195806161569131
$ws.Range("D6").Value = 1.789611489588995
$ws.Range("E6").Value = 1.789611489588995
$ws.Range("F6").Value = 1.411539595919041
$ws.Range("G6").Value = 1.411539595919041
$ws.Range("H6").Value = 0.4225844915764418
$ws.Range("I6").Value = 0.9337978568030243
$ws.Range("J6").Value = 0.8594090520211481
$ws.Range("K6").Value = 1.411539595919041
$ws.Range("L6").Value = 0.8195806161569131
$ws.Range("M6").Value = 1.304596052872954
$ws.Range("N6").Value = 1.304596052872954
$ws.Range("O6").Value = 1.01059219910745
$ws.Range("P6").Value = 1.340243900554983
$ws.Range("Q6").Value = 1.340243900554983
$ws.Range("R6").Value = 1.358067824395998
$ws.Range("S6").Value = 1.358067824395998
$ws.Range("T6").Value = 1.039420517010927

# row 7: Spiral2.5
$ws.Range("C7").Value = 1.002116304104381
$ws.Range("D7").Value = 0.9992340446463303
$ws.Range("E7").Value = 0.9992340446463303
$ws.Range("F7").Value = 1.003667642808365
$ws.Range("G7").Value = 1.003667642808365
$ws.Range("H7").Value = 0.9979958538094204
$ws.Range("I7").Value = 1.001241138171986
$ws.Range("J7").Value = 0.9988069448159292
$ws.Range("K7").Value = 1.003667642808365
$ws.Range("L7").Value = 1.002116304104381
$ws.Range("M7").Value = 1.000675174375356
$ws.Range("N7").Value = 1.000675174375356
$ws.Range("O7").Value = 0.999782067520044
$ws.Range("P7").Value = 1.001672663853025
$ws.Range("Q7").Value = 1.001672663853025
$ws.Range("R7").Value = 1.00217140859186
$ws.Range("S7").Value = 1.00217140859186
$ws.Range("T7").Value = 1.000510321392735

# row 8: Spiral5
$ws.Range("C8").Value = 1.002507699741749
$ws.Range("D8").Value = 0.9924400756700537
$ws.Range("E8").Value = 0.9924400756700537
$ws.Range("F8").Value = 1.015217704962765
$ws.Range("G8").Value = 1.015217704962765
$ws.Range("H8").Value = 0.9952985499878191
$ws.Range("I8").Value = 0.9994254430496223
$ws.Range("J8").Value = 0.9991375993798122
$ws.Range("K8").Value = 1.015217704962765
$ws.Range("L8").Value = 1.002507699741749
$ws.Range("M8").Value = 0.9974738877059015
$ws.Range("N8").Value = 0.9974738877059015
$ws.Range("O8").Value = 0.9967487751332073
$ws.Range("P8").Value = 1.003388493458189
$ws.Range("Q8").Value = 1.003388493458189
$ws.Range("R8").Value = 1.006345796334333
$ws.Range("S8").Value = 1.006345796334333
$ws.Range("T8").Value = 1.000671178798637

# row 9: Spiral7.5
$ws.Range("C9").Value = 1.017033413325553
$ws.Range("D9").Value = 0.9930889611819136
$ws.Range("E9").Value = 0.9930889611819136
$ws.Range("F9").Value = 1.062552636249771
$ws.Range("G9").Value = 1.062552636249771
$ws.Range("H9").Value = 0.9930936692391228
$ws.Range("I9").Value = 1.013065535706044
$ws.Range("J9").Value = 0.9918019649946318
$ws.Range("K9").Value = 1.062552636249771
$ws.Range("L9").Value = 1.017033413325553
$ws.Range("M9").Value = 1.005061187253733
$ws.Range("N9").Value = 1.005061187253733
$ws.Range("O9").Value = 1.001072014582197
$ws.Range("P9").Value = 1.024225003585746
$ws.Range("Q9").Value = 1.024225003585746
$ws.Range("R9").Value = 1.033806911751752
$ws.Range("S9").Value = 1.033806911751752
$ws.Range("T9").Value = 1.011772696782839

# row 10: Spiral10
$ws.Range("C10").Value = 1.009463378917367
$ws.Range("D10").Value = 0.9827206559784334
$ws.Range("E10").Value = 0.9827206559784334
$ws.Range("F10").Value = 1.074543995282749
$ws.Range("G10").Value = 1.074543995282749
$ws.Range("H10").Value = 0.9826155931610463
$ws.Range("I10").Value = 1.013911826369526
$ws.Range("J10").Value = 0.9970198950437998
$ws.Range("K10").Value = 1.074543995282749
$ws.Range("L10").Value = 1.009463378917367
$ws.Range("M10").Value = 0.9960920174479
$ws.Range("N10").Value = 0.9960920174479
$ws.Range("O10").Value = 0.9915998760189488
$ws.Range("P10").Value = 1.022242676726183
$ws.Range("Q10").Value = 1.022242676726183
$ws.Range("R10").Value = 1.035318006365324
$ws.Range("S10").Value = 1.035318006365324
$ws.Range("T10").Value = 1.010045890792153

# row 11: Spiral15
$ws.Range("C11").Value = 1.029046733070148
$ws.Range("D11").Value = 1.073930944309927
$ws.Range("E11").Value = 1.073930944309927
$ws.Range("F11").Value = 1.299474744863165
$ws.Range("G11").Value = 1.299474744863165
$ws.Range("H11").Value = 1.0620795190162
$ws.Range("I11").Value = 0.970119260442948
$ws.Range("J11").Value = 0.944353215695659
$ws.Range("K11").Value = 1.299474744863165
$ws.Range("L11").Value = 1.029046733070148
$ws.Range("M11").Value = 1.051488838690037
$ws.Range("N11").Value = 1.051488838690037
$ws.Range("O11").Value = 1.055019065465425
$ws.Range("P11").Value = 1.134150807414413
$ws.Range("Q11").Value = 1.134150807414413
$ws.Range("R11").Value = 1.175481791776601
$ws.Range("S11").Value = 1.175481791776601
$ws.Range("T11").Value = 1.063167402899674

# row 12: OffsetF45
$ws.Range("C12").Value = 1.200527794721968
$ws.Range("D12").Value = 0.02839472598084422
$ws.Range("E12").Value = 0.02839472598084422
$ws.Range("F12").Value = 0.8442444508594681
$ws.Range("G12").Value = 0.8442444508594681
$ws.Range("H12").Value = 1.667427177532818
$ws.Range("I12").Value = 0.007220350187736538
$ws.Range("J12").Value = 1.46987205148518
$ws.Range("K12").Value = 0.8442444508594681
$ws.Range("L12").Value = 1.200527794721968
$ws.Range("M12").Value = 0.6144612603514059
$ws.Range("N12").Value = 0.6144612603514059
$ws.Range("O12").Value = 0.9654498994118766
$ws.Range("P12").Value = 0.6910556571874267
$ws.Range("Q12").Value = 0.6910556571874267
$ws.Range("R12").Value = 0.7293528556054371
$ws.Range("S12").Value = 0.7293528556054371
$ws.Range("T12").Value = 0.8696144251280025

# row 13: OffsetA45
$ws.Range("C13").Value = 0.4933474806868863
$ws.Range("D13").Value = 3.312680609107422
$ws.Range("E13").Value = 3.312680609107422
$ws.Range("F13").Value = 1.098291033998245
$ws.Range("G13").Value = 1.098291033998245
$ws.Range("H13").Value = 0.6057433515298005
$ws.Range("I13").Value = 1.629499496735312
$ws.Range("J13").Value = 0.6647097133222074
$ws.Range("K13").Value = 1.098291033998245
$ws.Range("L13").Value = 0.4933474806868863
$ws.Range("M13").Value = 1.903014044897154
$ws.Range("N13").Value = 1.903014044897154
$ws.Range("O13").Value = 1.47059048044137
$ws.Range("P13").Value = 1.634773041264184
$ws.Range("Q13").Value = 1.634773041264184
$ws.Range("R13").Value = 1.500652539447699
$ws.Range("S13").Value = 1.500652539447699
$ws.Range("T13").Value = 1.300711947563312

# row 14: OffsetFTD
$ws.Range("C14").Value = 0.1746352361111097
$ws.Range("D14").Value = 3.41978210220953
$ws.Range("E14").Value = 3.41978210220953
$ws.Range("F14").Value = 1.919705994951049
$ws.Range("G14").Value = 1.919705994951049
$ws.Range("H14").Value = 1.834740474472013
$ws.Range("I14").Value = 0.1456247915945228
$ws.Range("J14").Value = 1.059089774778224
$ws.Range("K14").Value = 1.919705994951049
$ws.Range("L14").Value = 0.1746352361111097
$ws.Range("M14").Value = 1.79720866916032
$ws.Range("N14").Value = 1.79720866916032
$ws.Range("O14").Value = 1.809719270930884
$ws.Range("P14").Value = 1.838041111090563
$ws.Range("Q14").Value = 1.838041111090563
$ws.Range("R14").Value = 1.858457332055684
$ws.Range("S14").Value = 1.858457332055684
$ws.Range("T14").Value = 1.425596395686075

# row 15: OffsetATD
$ws.Range("C15").Value = 1.232545222439271
$ws.Range("D15").Value = 0.002627809240744217
$ws.Range("E15").Value = 0.002627809240744217
$ws.Range("F15").Value = 0.02078027073299304
$ws.Range("G15").Value = 0.02078027073299304
$ws.Range("H15").Value = 1.68803231914856
$ws.Range("I15").Value = 1.369476199899434
$ws.Range("J15").Value = 0.9313387374426676
$ws.Range("K15").Value = 0.02078027073299304
$ws.Range("L15").Value = 1.232545222439271
$ws.Range("M15").Value = 0.6175865158400077
$ws.Range("N15").Value = 0.6175865158400077
$ws.Range("O15").Value = 0.9744017836095252
$ws.Range("P15").Value = 0.4186511008043361
$ws.Range("Q15").Value = 0.4186511008043362
$ws.Range("R15").Value = 0.3191833932865004
$ws.Range("S15").Value = 0.3191833932865004
$ws.Range("T15").Value = 0.874133426483945

# ---- Step 5: populate the 4 new Holden rows (16-19) ----
# row 16: Holden2.5
$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "Holden2.5"
$ws.Range("C16").Value = 0.8723190885958978
$ws.Range("D16").Value = 0.1982716240525943
$ws.Range("E16").Value = 0.1982716240525943
$ws.Range("F16").Value = 1.755008585942779
$ws.Range("G16").Value = 1.755008585942779
$ws.Range("H16").Value = 0.2163825151528899
$ws.Range("I16").Value = 0.6855133870063155
$ws.Range("J16").Value = 1.022125262570573
$ws.Range("K16").Value = 1.755008585942779
$ws.Range("L16").Value = 0.8723190885958978
$ws.Range("M16").Value = 0.5352953563242461
$ws.Range("N16").Value = 0.5352953563242461
$ws.Range("O16").Value = 0.428991075933794
$ws.Range("P16").Value = 0.9418664328637569
$ws.Range("Q16").Value = 0.9418664328637569
$ws.Range("R16").Value = 1.145151971133512
$ws.Range("S16").Value = 1.145151971133512
$ws.Range("T16").Value = 0.7916034105535082

# row 17: Holden5
$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "Holden5"
$ws.Range("C17").Value = 0.9850677104148882
$ws.Range("D17").Value = 0.3963049241687018
$ws.Range("E17").Value = 0.3963049241687018
$ws.Range("F17").Value = 1.475237453348569
$ws.Range("G17").Value = 1.475237453348569
$ws.Range("H17").Value = 0.4187003194036331
$ws.Range("I17").Value = 0.9226522745937056
$ws.Range("J17").Value = 1.016035199587499
$ws.Range("K17").Value = 1.475237453348569
$ws.Range("L17").Value = 0.9850677104148882
$ws.Range("M17").Value = 0.690686317291795
$ws.Range("N17").Value = 0.690686317291795
$ws.Range("O17").Value = 0.6000243179957411
$ws.Range("P17").Value = 0.952203362644053
$ws.Range("Q17").Value = 0.952203362644053
$ws.Range("R17").Value = 1.082961885320182
$ws.Range("S17").Value = 1.082961885320182
$ws.Range("T17").Value = 0.8689996469194995

# row 18: Holden10
$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "Holden10"
$ws.Range("C18").Value = 1.212901290104071
$ws.Range("D18").Value = 0.7845849941778587
$ws.Range("E18").Value = 0.7845849941778587
$ws.Range("F18").Value = 0.9530112008586459
$ws.Range("G18").Value = 0.9530112008586459
$ws.Range("H18").Value = 0.8212349010279795
$ws.Range("I18").Value = 1.395436807946298
$ws.Range("J18").Value = 1.012784851273102
$ws.Range("K18").Value = 0.9530112008586459
$ws.Range("L18").Value = 1.212901290104071
$ws.Range("M18").Value = 0.9987431421409646
$ws.Range("N18").Value = 0.9987431421409646
$ws.Range("O18").Value = 0.9395737284366362
$ws.Range("P18").Value = 0.983499161713525
$ws.Range("Q18").Value = 0.983499161713525
$ws.Range("R18").Value = 0.9758771714998052
$ws.Range("S18").Value = 0.9758771714998052
$ws.Range("T18").Value = 1.029992340897992

# row 19: Holden15
$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "Holden15"
$ws.Range("C19").Value = 1.346443177847401
$ws.Range("D19").Value = 0.7590534233716273
$ws.Range("E19").Value = 0.7590534233716273
$ws.Range("F19").Value = 0.8592724195351429
$ws.Range("G19").Value = 0.8592724195351429
$ws.Range("H19").Value = 0.8390952770782442
$ws.Range("I19").Value = 1.400832096336024
$ws.Range("J19").Value = 0.9456020791841662
$ws.Range("K19").Value = 0.8592724195351429
$ws.Range("L19").Value = 1.346443177847401
$ws.Range("M19").Value = 1.052748300609514
$ws.Range("N19").Value = 1.052748300609514
$ws.Range("O19").Value = 0.9815306260990909
$ws.Range("P19").Value = 0.9882563402513904
$ws.Range("Q19").Value = 0.9882563402513904
$ws.Range("R19").Value = 0.9560103600723285
$ws.Range("S19").Value = 0.9560103600723285
$ws.Range("T19").Value = 1.025049745558768

# ---- Step 6: refresh the (now shifted) HexGrid rows 20-23, including the index column ----
# row 20: HexGrid-90degTilt2.5degRes
$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "HexGrid-90degTilt2.5degRes"
$ws.Range("C20").Value = 1.008975085568029
$ws.Range("D20").Value = 1.001408900416044
$ws.Range("E20").Value = 1.001408900416044
$ws.Range("F20").Value = 0.9889247942658974
$ws.Range("G20").Value = 0.9889247942658974
$ws.Range("H20").Value = 1.001140825011323
$ws.Range("I20").Value = 1.005531437017992
$ws.Range("J20").Value = 0.9983746269812633
$ws.Range("K20").Value = 0.9889247942658974
$ws.Range("L20").Value = 1.008975085568029
$ws.Range("M20").Value = 1.005191992992037
$ws.Range("N20").Value = 1.005191992992037
$ws.Range("O20").Value = 1.003841603665132
$ws.Range("P20").Value = 0.9997695934166568
$ws.Range("Q20").Value = 0.9997695934166568
$ws.Range("R20").Value = 0.997058393628967
$ws.Range("S20").Value = 0.997058393628967
$ws.Range("T20").Value = 1.000725944876758

# row 21: HexGrid-90degTilt5degRes
$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "HexGrid-90degTilt5degRes"
$ws.Range("C21").Value = 0.9924101642594766
$ws.Range("D21").Value = 1.000319665283988
$ws.Range("E21").Value = 1.000319665283988
$ws.Range("F21").Value = 1.022235747652223
$ws.Range("G21").Value = 1.022235747652223
$ws.Range("H21").Value = 0.9981719549974212
$ws.Range("I21").Value = 0.9991074780891017
$ws.Range("J21").Value = 1.000587706943871
$ws.Range("K21").Value = 1.022235747652223
$ws.Range("L21").Value = 0.9924101642594766
$ws.Range("M21").Value = 0.9963649147717324
$ws.Range("N21").Value = 0.9963649147717324
$ws.Range("O21").Value = 0.9969672615136287
$ws.Range("P21").Value = 1.004988525731896
$ws.Range("Q21").Value = 1.004988525731896
$ws.Range("R21").Value = 1.009300331211977
$ws.Range("S21").Value = 1.009300331211977
$ws.Range("T21").Value = 1.002138786204347

# row 22: HexGrid-90degTilt10degRes
$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "HexGrid-90degTilt10degRes"
$ws.Range("C22").Value = 1.008637648358219
$ws.Range("D22").Value = 0.9909691756832414
$ws.Range("E22").Value = 0.9909691756832414
$ws.Range("F22").Value = 1.111972894212945
$ws.Range("G22").Value = 1.111972894212945
$ws.Range("H22").Value = 0.9914177580555773
$ws.Range("I22").Value = 1.018284271939523
$ws.Range("J22").Value = 1.004392640781794
$ws.Range("K22").Value = 1.111972894212945
$ws.Range("L22").Value = 1.008637648358219
$ws.Range("M22").Value = 0.9998034120207302
$ws.Range("N22").Value = 0.9998034120207302
$ws.Range("O22").Value = 0.9970081940323459
$ws.Range("P22").Value = 1.037193239418135
$ws.Range("Q22").Value = 1.037193239418135
$ws.Range("R22").Value = 1.055888153116837
$ws.Range("S22").Value = 1.055888153116837
$ws.Range("T22").Value = 1.020945731505217

# row 23: HexGrid-90degTilt15degRes
$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "HexGrid-90degTilt15degRes"
$ws.Range("C23").Value = 1.157169476054419
$ws.Range("D23").Value = 0.6669968676775767
$ws.Range("E23").Value = 0.6669968676775767
$ws.Range("F23").Value = 0.9605819309588933
$ws.Range("G23").Value = 0.9605819309588933
$ws.Range("H23").Value = 1.074046226701863
$ws.Range("I23").Value = 0.9750234546700978
$ws.Range("J23").Value = 0.9302875081231422
$ws.Range("K23").Value = 0.9605819309588933
$ws.Range("L23").Value = 1.157169476054419
$ws.Range("M23").Value = 0.9120831718659979
$ws.Range("N23").Value = 0.9120831718659979
$ws.Range("O23").Value = 0.9660708568112861
$ws.Range("P23").Value = 0.928249424896963
$ws.Range("Q23").Value = 0.928249424896963
$ws.Range("R23").Value = 0.9363325514124456
$ws.Range("S23").Value = 0.9363325514124456
$ws.Range("T23").Value = 0.9606842440309987

